$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# -- Simple value updates (Metadata sheet) --
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# -- Give row 16 the same formatting as the rest of the data rows --
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# -- Shift Description/Purpose/Copyright/Immutable down one row (12-15 -> 13-16) --
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").ClearContents()

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").ClearContents()

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Platelet morphology panel - Blood (58406-0)"

# -- Row 12 becomes the new "Jurisdiction" row (keeps the pre-existing data-row style) --
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
